# Auto-generated Excel COM-interop script to update cryptos list data
# (mirrors the source XML diff: price / volume updates and an ONDO <-> EnergySwap row swap)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = '66.476.79'
$ws.Range("E2").Value = '  -0.28%  '
$ws.Range("D3").Value = '3.509.15'
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '604.51'
$ws.Range("E5").Value = '  -1.02%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.44'
$ws.Range("E6").Value = '  -4.59%  '
$ws.Range("D7").Value = '3.509.94'
$ws.Range("E7").Value = '  -3.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  -0.15%  '
$ws.Range("E9").Value = '  +3.28%  '
$ws.Range("E10").Value = '  -3.11%  '
$ws.Range("E11").Value = '  -5.54%  '
$ws.Range("E12").Value = '  -3.29%  '
$ws.Range("D13").Value = '4.096.58'
$ws.Range("E13").Value = '  -3.27%  '
$ws.Range("E14").Value = '  -7.38%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '28.63'
$ws.Range("E15").Value = '  -4.51%  '
$ws.Range("D16").Value = '3.519.69'
$ws.Range("E16").Value = '  -2.72%  '
$ws.Range("E17").Value = '  -0.24%  '
$ws.Range("D18").Value = '66.334.87'
$ws.Range("E18").Value = '  -0.63%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.63'
$ws.Range("E19").Value = '  -8.93%  '
$ws.Range("E20").Value = '  -4.56%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.56'
$ws.Range("E21").Value = '  -3.62%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '421.19'
$ws.Range("E22").Value = '  -1.68%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.586'
$ws.Range("E23").Value = '  -5.46%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '76.66'
$ws.Range("E24").Value = '  -2.79%  '
$ws.Range("D25").Value = '3.652.93'
$ws.Range("E25").Value = '  -2.99%  '
$ws.Range("E26").Value = '  -0.09%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000113'
$ws.Range("E27").Value = '  -8.87%  '
$ws.Range("E28").Value = '  -3.26%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.76'
$ws.Range("E29").Value = '  -7.81%  '
$ws.Range("E30").Value = '  -7.27%  '
$ws.Range("E31").Value = '  +0.05%  '
$ws.Range("D32").Value = '3.513.06'
$ws.Range("E32").Value = '  -2.94%  '
$ws.Range("E33").Value = '  -4.00%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '24.12'
$ws.Range("E34").Value = '  -5.29%  '
$ws.Range("E35").Value = '  +0.00%  '
$ws.Range("E36").Value = '  -10.54%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '7.49'
$ws.Range("E37").Value = '  -5.21%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.63'
$ws.Range("E38").Value = '  -4.90%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '173.42'
$ws.Range("E39").Value = '  -2.27%  '
$ws.Range("E40").Value = '  -8.62%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0803'
$ws.Range("E41").Value = '  -6.92%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.93'
$ws.Range("E42").Value = '  -5.86%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.850'
$ws.Range("E43").Value = '  -5.69%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '45.50'
$ws.Range("E44").Value = '  -1.60%  '
$ws.Range("E45").Value = '  -7.43%  '
$ws.Range("E46").Value = '  +0.10%  '
$ws.Range("E47").Value = '  -10.45%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.02'
$ws.Range("E48").Value = '  -2.66%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '22.93'
$ws.Range("E49").Value = '  -4.25%  '
$ws.Range("B50").Value = 'ONDO'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.11'
$ws.Range("E50").Value = '  -5.10%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.894'
$ws.Range("E51").Value = '  -7.43%  '
